# Update "想去人数" (number of people interested) figures in the
# 北京-漫展信息 workbook, as produced by the gh-pages data refresh at 456a3b4.
#
# The same underlying event data appears on multiple sheets:
#   展览    (sheet1, "Exhibitions")
#   演出    (sheet2, "Performances")
#   本地生活 (sheet3, "Local life")
#   全部类型 (sheet4, "All types" - combines the above)
# so matching rows must be updated on every sheet that contains them.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# ---- 展览 (sheet1) ----
$ws1.Range("F3").Value = 2860
$ws1.Range("F5").Value = 6386
$ws1.Range("F6").Value = 2483
$ws1.Range("F10").Value = 2897
$ws1.Range("F11").Value = 45
$ws1.Range("F13").Value = 7140
$ws1.Range("F14").Value = 301
$ws1.Range("F15").Value = 21
$ws1.Range("F18").Value = 113
$ws1.Range("F20").Value = 8407
$ws1.Range("F21").Value = 18
$ws1.Range("F28").Value = 79
$ws1.Range("F30").Value = 15
$ws1.Range("F32").Value = 57
$ws1.Range("F33").Value = 78
$ws1.Range("F34").Value = 2600
$ws1.Range("F37").Value = 33
$ws1.Range("F38").Value = 1166
$ws1.Range("F40").Value = 684
$ws1.Range("F41").Value = 3678
$ws1.Range("F42").Value = 5
$ws1.Range("F43").Value = 181
$ws1.Range("F44").Value = 18
$ws1.Range("F45").Value = 1189
$ws1.Range("F46").Value = 176
$ws1.Range("F47").Value = 20

# ---- 演出 (sheet2) ----
$ws2.Range("F6").Value = 4

# ---- 全部类型 (sheet4) ----
$ws4.Range("F3").Value = 2860
$ws4.Range("F6").Value = 6386
$ws4.Range("F7").Value = 2483
$ws4.Range("F12").Value = 2897
$ws4.Range("F13").Value = 45
$ws4.Range("F17").Value = 7140
$ws4.Range("F18").Value = 301
$ws4.Range("F21").Value = 113
$ws4.Range("F23").Value = 8407
$ws4.Range("F29").Value = 79
$ws4.Range("F32").Value = 57
$ws4.Range("F34").Value = 78
$ws4.Range("F35").Value = 2600
$ws4.Range("F38").Value = 33
$ws4.Range("F39").Value = 1166
$ws4.Range("F40").Value = 684
$ws4.Range("F42").Value = 3678
$ws4.Range("F43").Value = 181
$ws4.Range("F44").Value = 18
$ws4.Range("F46").Value = 1189
$ws4.Range("F47").Value = 176
$ws4.Range("F48").Value = 20
